# Rename the seven "*img" sheets to "img*" (keeps the trailing letter but
# moves the "img" prefix to the front of the sheet name).
$wb = $excel.ActiveWorkbook

$renames = @{
    "himg" = "imgh"
    "timg" = "imgt"
    "simg" = "imgs"
    "gimg" = "imgg"
    "wimg" = "imgw"
    "bimg" = "imgb"
    "eimg" = "imge"
}

foreach ($oldName in $renames.Keys) {
    $sheet = $wb.Worksheets.Item($oldName)
    $sheet.Name = $renames[$oldName]
}

# The last sheet (formerly "eimg", now "imge") becomes the active/selected tab.
$wb.Worksheets.Item("imge").Activate()
